$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

# Original paragraph is split across three runs: "Below", " ", "section-level".
# Remove the single-space run, then fold the leading word back together with
# the trailing space so the whole paragraph collapses into one run reading
# "Below section-level" (matching the target OOXML: a single <a:r> with an
# empty <a:rPr/>).
$space = $tr.Characters(6, 1)
$space.Text = ""

$lead = $tr.Characters(1, 5)
$lead.Text = "Below "

$whole = $tr.Characters(1, $tr.Length)
$whole.Text = "Below section-level"
